# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B describes the dimension/measure metadata for "municipio-nombre".
# It is being re-curated from an iaest-measure to an sdmx-dimension, and its
# type/label rows are updated to match the new curated dimension.
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("B3").Value = "dim"
$ws.Range("B4").Value = "URI-Municipio"
